$wb = $excel.ActiveWorkbook

# --- Rename sheets and drop the now-unused third sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "PID"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Temp"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Delete() | Out-Null

# --- Rebuild sheet "PID" (was Hoja1) with the new parameter layout ---
$ws1.Cells.Clear() | Out-Null

$ws1.Range("B2").Value = 5400
$ws1.Range("D2").Value = 7200

$ws1.Range("B3").Value = 200
$ws1.Range("D3").Value = 150

$ws1.Range("B4").Value = 2400
$ws1.Range("D4").Value = 1600

$ws1.Range("B5").Value = 1
$ws1.Range("D5").Value = 1

$ws1.Range("B6").Value = 1000
$ws1.Range("D6").Value = 1000

$ws1.Range("B7").Value = 0
$ws1.Range("D7").Value = 0

$ws1.Range("B8").Value = 20
$ws1.Range("D8").Value = 20

$ws1.Range("B9").Value = 795
$ws1.Range("D9").Value = 795

$ws1.Range("B10").Value = 3000
$ws1.Range("D10").Value = 3000

$ws1.Range("C16").Select() | Out-Null

# --- Rebuild sheet "Temp" (was Hoja2) with its new data ---
$ws2.Cells.Clear() | Out-Null

$ws2.Range("B2").Value = 704
$ws2.Range("D2").Value = 9.94

$ws2.Range("B3").Value = 1544
$ws2.Range("D3").Value = 205

$ws2.Range("C5").Value = 5

$ws2.Range("B10").Select() | Out-Null

$ws1.Activate() | Out-Null
